$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 6: switch the table's style from the deck's custom
#    "Table_0" style ({FBE09CBE-9BDB-4000-B46A-FF644B45D47B}) to the built-in
#    PowerPoint table style {53097EEF-DF33-45F5-B609-1E793796033B}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{53097EEF-DF33-45F5-B609-1E793796033B}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colours: the deck's theme ("Integral") is re-coloured to match the
#    stock "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# ---------------------------------------------------------------------------
function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeSlide = $p.Slides.Item(1)
$colorScheme = $themeSlide.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToVbaRgb($officeColors[$i - 1])
}
